{"js": "// Split the bullet \"CSPO Statements. Populated Kinds.\" (ilvl=2, numId=1)\n// into three separate bullets at the same list level:\n//   1. \"CSPO Statements.\"\n//   2. \"Populated Kinds.\"\n//   3. \"(Context, Stat1Subject, Stat2Property, Stat3Object); (when Kinds Alignment matches).\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst ORIGINAL_TEXT = \"CSPO Statements. Populated Kinds.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === ORIGINAL_TEXT) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\n    'Could not find the paragraph with text \"' + ORIGINAL_TEXT + '\"'\n  );\n}\n\n// Trim the original paragraph's text down to just the first sentence.\ntarget.getRange().insertText(\"CSPO Statements.\", \"Replace\");\n\n// Insert the two new bullets right after it. `insertParagraph` on a\n// paragraph copies that paragraph's formatting (numbering/indent/shading),\n// so the new bullets stay at the same list level as the original.\nconst secondParagraph = target.insertParagraph(\"Populated Kinds.\", \"After\");\nsecondParagraph.insertParagraph(\n  \"(Context, Stat1Subject, Stat2Property, Stat3Object); (when Kinds Alignment matches).\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Split the bullet \"CSPO Statements. Populated Kinds.\" (ilvl=2, numId=1)\n# into three separate bullets at the same list level:\n#   1. \"CSPO Statements.\"\n#   2. \"Populated Kinds.\"\n#   3. \"(Context, Stat1Subject, Stat2Property, Stat3Object); (when Kinds Alignment matches).\"\n$d = $word.ActiveDocument\n\n$originalText = \"CSPO Statements. Populated Kinds.\"\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (chr 13).\n    if ($p.Range.Text -eq ($originalText + [char]13)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the paragraph with text '$originalText'\"\n}\n\n# Trim the original paragraph's text down to just the first sentence\n# (assigning Range.Text only replaces the text, keeping the paragraph mark\n# and its formatting/list membership intact).\n$target.Range.Text = \"CSPO Statements.\"\n\n# Insert a new paragraph after it, inheriting the same paragraph formatting\n# (numbering/indent/shading), then set its text.\n$target.Range.InsertParagraphAfter()\n$secondParagraph = $target.Next()\n$secondParagraph.Range.Text = \"Populated Kinds.\"\n\n# Insert the third paragraph the same way.\n$secondParagraph.Range.InsertParagraphAfter()\n$thirdParagraph = $secondParagraph.Next()\n$thirdParagraph.Range.Text = \"(Context, Stat1Subject, Stat2Property, Stat3Object); (when Kinds Alignment matches).\"\n"}
